# Applies the "more intrinsics, execution tweeks, unit test love for lexer"
# update to the Supported.xlsx workbook: marks several additional
# functions/statements as supported (Token/Parse/Eval columns) and adds a
# few explanatory notes in the "Unsupported operations" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows that gain full Token/Parse/Eval ("X") support.
$newlySupported = @(55, 61, 63, 68, 69, 75, 88, 89, 94, 102, 143, 145)

foreach ($r in $newlySupported) {
    $ws.Cells.Item($r, 2).Value = "X"
    $ws.Cells.Item($r, 3).Value = "X"
    $ws.Cells.Item($r, 4).Value = "X"
}

# Explanatory notes in column E ("Unsupported operations").
$ws.Cells.Item(64, 5).Value = "No way to implement"
$ws.Cells.Item(89, 5).Value = "returns 0 until printing is working"
$ws.Cells.Item(94, 5).Value = "function, not statement"

# Move the selection to where the author left off editing.
$ws.Range("B149").Select()
